$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlaceEntity")

# Dialogue Typing Text Add:
# Lorena's face emote on her spoken lines (column C, rows where she is the
# speaker) changes from "Lorena_Happy" to "Lorena_Sad". Re-writing the cell
# from scratch (Clear + Value) drops the inherited table-border style so the
# cell reverts to the default/no explicit style, matching the edit.
$rows = @(15, 17, 18, 20, 21)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Clear() | Out-Null
    $cell.Value = "Lorena_Sad"
}

# Move the active selection from D15 to C15.
$ws.Range("C15").Select() | Out-Null
